$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "organizations:" block - append ", church" to the hospital/.../
#    prison paragraph, then insert the new "not organizations:" /
#    "hotel..." / "olympics..., program" paragraphs *before* the
#    hidden "_GoBack" bookmark that sits at the end of that
#    paragraph, and finish with the brand-new "But proper names..."
#    sentence (which is where the bookmark ends up after the edit).
#    We do this as a single Find/Replace so the bookmark (invisible
#    to this runtime's Bookmarks collection) keeps tracking the
#    original insertion point, landing in its own (temporarily empty)
#    paragraph that we then fill in via a direct Range.Text set -
#    that keeps any new text *before* the bookmark, matching Word's
#    own behaviour.
# ------------------------------------------------------------------

$d.Content.Find.Execute(
    ", prison", $true, $false, $false, $false, $false, $true, 1, $false,
    ", prison, church" + "`r" +
    "not organizations:" + "`r" +
    "  hotel, motel, airbase, office, district, airport, resort, marathon," + "`r" +
    "  olympics, shopping mall, auditorium, program" + "`r",
    2)

# Locate the freshly-split, still-empty paragraph that now holds the
# bookmark (it directly follows the "...auditorium, program" one) and
# give it the final sentence - this inserts the text ahead of the
# bookmark, exactly like the diff wants.
$bookmarkPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd() -eq "  olympics, shopping mall, auditorium, program") {
        $bookmarkPara = $d.Paragraphs.Item($i + 1)
        break
    }
}
$bookmarkPara.Range.Text = "But proper names of airbases, police departments, etc. are organizations"
$bookmarkPara.Range.Font.Name = "Courier New"
$bookmarkPara.Range.Font.NameAscii = "Courier New"
$bookmarkPara.Range.Font.NameBi = "Courier New"

# The old copies of "not organizations:" / "hotel, motel, ..." /
# "olympics, shopping mall, auditorium" (without ", program") are now
# duplicated right after the new "But proper names..." paragraph -
# delete those three paragraphs (the following, originally blank,
# paragraph is left untouched).
$oldNotOrgIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd() -eq "not organizations:") {
        $oldNotOrgIdx = $i
    }
}
$delStart = $d.Paragraphs.Item($oldNotOrgIdx).Range.Start
$delEnd = $d.Paragraphs.Item($oldNotOrgIdx + 2).Range.End
$d.Range($delStart, $delEnd).Delete()

# ------------------------------------------------------------------
# 2) "insurgents" gains a new ", community" run.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "  insurgents", $true, $false, $false, $false, $false, $true, 1, $false,
    "  insurgents, community", 2)

# ------------------------------------------------------------------
# 3) "  police, navy," loses its trailing ",".
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    ", navy,", $true, $false, $false, $false, $false, $true, 1, $false,
    ", navy", 2)

# ------------------------------------------------------------------
# 4) New paragraphs at the very end of the document: a blank line,
#    "School: can modify a division within the school", an
#    entity/slotfill line, and a trailing blank paragraph.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$tail = $lastPara.Range.Duplicate
$tail.Collapse(0)
$tail.InsertAfter(
    "`r`r" +
    "School: can modify a division within the school" + "`r" +
    "  entity: " + "Johns Hopkins School of Medicine" + "    slotfill: " + "Johns Hopkins" + "`r")
